# Auto-generated Excel COM-interop script applying the 2026-02-07 01:49
# meteocat daily-summary refresh (extraction timestamps + updated readings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a bare "NN%" need the cell pre-formatted as Text,
# otherwise Excel (like the real app) auto-converts "54%" into the number 0.54
# with a Percentage number format instead of keeping literal text.
$ws.Range('H4').NumberFormat = "@"
$ws.Range('H5').NumberFormat = "@"
$ws.Range('H7').NumberFormat = "@"
$ws.Range('H8').NumberFormat = "@"
$ws.Range('H9').NumberFormat = "@"
$ws.Range('H11').NumberFormat = "@"
$ws.Range('H12').NumberFormat = "@"
$ws.Range('H15').NumberFormat = "@"
$ws.Range('H16').NumberFormat = "@"
$ws.Range('H20').NumberFormat = "@"
$ws.Range('H21').NumberFormat = "@"
$ws.Range('H22').NumberFormat = "@"
$ws.Range('H24').NumberFormat = "@"
$ws.Range('H25').NumberFormat = "@"
$ws.Range('H26').NumberFormat = "@"
$ws.Range('H27').NumberFormat = "@"
$ws.Range('H29').NumberFormat = "@"
$ws.Range('H30').NumberFormat = "@"
$ws.Range('H32').NumberFormat = "@"
$ws.Range('H33').NumberFormat = "@"
$ws.Range('H34').NumberFormat = "@"

$ws.Range('E2').Value = '2026-02-07 01:47:57'
$ws.Range('N2').Value = '-1.2 °C 1:29 TU'
$ws.Range('O2').Value = '-1.0 °C'
$ws.Range('E3').Value = '2026-02-07 01:48:00'
$ws.Range('N3').Value = '-5.2 °C 1:29 TU'
$ws.Range('O3').Value = '-4.4 °C'
$ws.Range('E4').Value = '2026-02-07 01:48:02'
$ws.Range('H4').Value = '54%'
$ws.Range('N4').Value = '9.9 °C 1:26 TU'
$ws.Range('O4').Value = '12.1 °C'
$ws.Range('E5').Value = '2026-02-07 01:48:05'
$ws.Range('H5').Value = '68%'
$ws.Range('L5').Value = '13.3 km/h - 216º 1:06 TU'
$ws.Range('O5').Value = '9.7 °C'
$ws.Range('E6').Value = '2026-02-07 01:48:08'
$ws.Range('J6').Value = '1001.9 hPa'
$ws.Range('N6').Value = '12.0 °C 1:29 TU'
$ws.Range('O6').Value = '12.5 °C'
$ws.Range('E7').Value = '2026-02-07 01:48:10'
$ws.Range('H7').Value = '71%'
$ws.Range('N7').Value = '7.8 °C 1:29 TU'
$ws.Range('O7').Value = '8.6 °C'
$ws.Range('E8').Value = '2026-02-07 01:48:13'
$ws.Range('H8').Value = '92%'
$ws.Range('M8').Value = '5.9 °C 1:09 TU'
$ws.Range('O8').Value = '5.3 °C'
$ws.Range('E9').Value = '2026-02-07 01:48:15'
$ws.Range('H9').Value = '99%'
$ws.Range('M9').Value = '5.3 °C 1:22 TU'
$ws.Range('O9').Value = '3.0 °C'
$ws.Range('E10').Value = '2026-02-07 01:48:18'
$ws.Range('M10').Value = '7.6 °C 1:29 TU'
$ws.Range('N10').Value = '6.4 °C 1:00 TU'
$ws.Range('O10').Value = '7.0 °C'
$ws.Range('E11').Value = '2026-02-07 01:48:21'
$ws.Range('H11').Value = '95%'
$ws.Range('I11').Value = '2.7 mm'
$ws.Range('O11').Value = '1.5 °C'
$ws.Range('E12').Value = '2026-02-07 01:48:23'
$ws.Range('H12').Value = '68%'
$ws.Range('L12').Value = '25.2 km/h - 266º 1:25 TU'
$ws.Range('N12').Value = '9.8 °C 1:13 TU'
$ws.Range('O12').Value = '10.7 °C'
$ws.Range('E13').Value = '2026-02-07 01:48:25'
$ws.Range('M13').Value = '8.1 °C 1:27 TU'
$ws.Range('N13').Value = '6.6 °C 1:06 TU'
$ws.Range('E14').Value = '2026-02-07 01:48:28'
$ws.Range('O14').Value = '-5.8 °C'
$ws.Range('E15').Value = '2026-02-07 01:48:31'
$ws.Range('H15').Value = '76%'
$ws.Range('O15').Value = '8.3 °C'
$ws.Range('E16').Value = '2026-02-07 01:48:33'
$ws.Range('H16').Value = '84%'
$ws.Range('N16').Value = '3.5 °C 1:29 TU'
$ws.Range('O16').Value = '4.2 °C'
$ws.Range('E17').Value = '2026-02-07 01:48:36'
$ws.Range('N17').Value = '3.4 °C 1:10 TU'
$ws.Range('O17').Value = '3.8 °C'
$ws.Range('E18').Value = '2026-02-07 01:48:39'
$ws.Range('N18').Value = '-6.2 °C 1:29 TU'
$ws.Range('O18').Value = '-6.1 °C'
$ws.Range('E19').Value = '2026-02-07 01:48:41'
$ws.Range('J19').Value = '1005.0 hPa'
$ws.Range('N19').Value = '5.1 °C 1:29 TU'
$ws.Range('E20').Value = '2026-02-07 01:48:44'
$ws.Range('H20').Value = '91%'
$ws.Range('L20').Value = '6.5 km/h - 239º 1:25 TU'
$ws.Range('N20').Value = '-4.3 °C 1:29 TU'
$ws.Range('E21').Value = '2026-02-07 01:48:46'
$ws.Range('H21').Value = '59%'
$ws.Range('N21').Value = '9.0 °C 1:21 TU'
$ws.Range('O21').Value = '10.3 °C'
$ws.Range('E22').Value = '2026-02-07 01:48:49'
$ws.Range('H22').Value = '94%'
$ws.Range('M22').Value = '6.8 °C 1:20 TU'
$ws.Range('O22').Value = '6.1 °C'
$ws.Range('E23').Value = '2026-02-07 01:48:51'
$ws.Range('L23').Value = '11.5 km/h - 38º 1:08 TU'
$ws.Range('M23').Value = '8.0 °C 1:10 TU'
$ws.Range('E24').Value = '2026-02-07 01:48:54'
$ws.Range('H24').Value = '78%'
$ws.Range('N24').Value = '10.2 °C 1:10 TU'
$ws.Range('O24').Value = '10.7 °C'
$ws.Range('E25').Value = '2026-02-07 01:48:56'
$ws.Range('H25').Value = '95%'
$ws.Range('I25').Value = '1.9 mm'
$ws.Range('J25').Value = '1004.8 hPa'
$ws.Range('N25').Value = '0.3 °C 1:26 TU'
$ws.Range('O25').Value = '1.1 °C'
$ws.Range('E26').Value = '2026-02-07 01:48:59'
$ws.Range('H26').Value = '80%'
$ws.Range('L26').Value = '18.4 km/h - 348º 1:27 TU'
$ws.Range('M26').Value = '-1.1 °C 1:26 TU'
$ws.Range('E27').Value = '2026-02-07 01:49:02'
$ws.Range('H27').Value = '96%'
$ws.Range('M27').Value = '9.9 °C 1:04 TU'
$ws.Range('N27').Value = '8.0 °C 1:29 TU'
$ws.Range('O27').Value = '8.5 °C'
$ws.Range('E28').Value = '2026-02-07 01:49:04'
$ws.Range('N28').Value = '4.1 °C 1:21 TU'
$ws.Range('O28').Value = '4.6 °C'
$ws.Range('E29').Value = '2026-02-07 01:49:07'
$ws.Range('H29').Value = '53%'
$ws.Range('N29').Value = '11.6 °C 1:25 TU'
$ws.Range('O29').Value = '12.5 °C'
$ws.Range('E30').Value = '2026-02-07 01:49:09'
$ws.Range('H30').Value = '78%'
$ws.Range('I30').Value = '0.3 mm'
$ws.Range('N30').Value = '-4.9 °C 1:28 TU'
$ws.Range('O30').Value = '-4.1 °C'
$ws.Range('E31').Value = '2026-02-07 01:49:12'
$ws.Range('N31').Value = '3.8 °C 1:29 TU'
$ws.Range('E32').Value = '2026-02-07 01:49:14'
$ws.Range('H32').Value = '65%'
$ws.Range('J32').Value = '1003.3 hPa'
$ws.Range('L32').Value = '28.4 km/h - 287º 1:14 TU'
$ws.Range('O32').Value = '11.5 °C'
$ws.Range('E33').Value = '2026-02-07 01:49:17'
$ws.Range('H33').Value = '95%'
$ws.Range('M33').Value = '8.1 °C 1:22 TU'
$ws.Range('O33').Value = '7.2 °C'
$ws.Range('E34').Value = '2026-02-07 01:49:19'
$ws.Range('H34').Value = '72%'
$ws.Range('N34').Value = '6.5 °C 1:17 TU'
$ws.Range('O34').Value = '7.5 °C'
$ws.Range('E35').Value = '2026-02-07 01:49:22'
$ws.Range('G35').Value = '201 cm'
$ws.Range('N35').Value = '-4.1 °C 1:01 TU'
$ws.Range('E36').Value = '2026-02-07 01:49:24'
$ws.Range('J36').Value = '1005.4 hPa'
